$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data added to the first table (rows 6-8): a "-g0 -Wl,--strip-all"
#    data row (row 7), a couple of extra header/annotation cells on row 6
#    (H6/I6: "w/alpine" / "w/debian-slim"), and the two new Docker-layer /
#    build-time cells in rows 6 and 8.
#    Write order matters: it controls the order in which brand-new strings
#    are appended to the shared-string table.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "17M"
$ws.Range("C7").Value = "57M"
$ws.Range("D7").Value = "41M"
$ws.Range("E7").Value = "22M"
$ws.Range("F7").Value = "321M"
$ws.Range("G7").Value = 0.011308576388888889
$ws.Range("G7").NumberFormat = "mm:ss.0"
$ws.Range("H7").Value = "546MB"

$ws.Range("H6").Value = "w/alpine"
$ws.Range("I6").Value = "w/debian-slim"
$ws.Range("F6").Value = "334M"

$ws.Range("G8").Value = 0.0077814120370370381
$ws.Range("G8").NumberFormat = "mm:ss.0"

# ---------------------------------------------------------------------------
# 2. Re-order the second table (rows 13-16) so the CFLAG rows follow the
#    same sequence as the first table (-O2, -Os, -g0, -Wl,--strip-all).
#    Capture the current values first, then rewrite the rows in the new
#    order. Row 16 (-Wl,--strip-all) does not move.
# ---------------------------------------------------------------------------
$row13 = @($ws.Range("A13").Value2, $ws.Range("B13").Value2, $ws.Range("C13").Value2, $ws.Range("D13").Value2, $ws.Range("E13").Value2)
$row14 = @($ws.Range("A14").Value2, $ws.Range("B14").Value2, $ws.Range("C14").Value2, $ws.Range("D14").Value2, $ws.Range("E14").Value2)
$row15 = @($ws.Range("A15").Value2, $ws.Range("B15").Value2, $ws.Range("C15").Value2, $ws.Range("D15").Value2, $ws.Range("E15").Value2)

# New row 13 = old row 14 ("-O2")
$ws.Range("A13").Value = $row14[0]
$ws.Range("B13").Value = $row14[1]
$ws.Range("C13").Value = $row14[2]
$ws.Range("D13").Value = $row14[3]
$ws.Range("E13").Value = $row14[4]

# New row 14 = old row 15 ("-Os"), no NumPy/SciPy figures
$ws.Range("A14").Value = $row15[0]
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = $row15[3]
$ws.Range("E14").Value = $row15[4]

# New row 15 = old row 13 ("-g0")
$ws.Range("A15").Value = $row13[0]
$ws.Range("B15").Value = $row13[1]
$ws.Range("C15").Value = $row13[2]
$ws.Range("D15").Value = $row13[3]
$ws.Range("E15").Value = $row13[4]

# Row 16 ("-Wl,--strip-all") is unchanged.

# Re-stamp column A's "quote prefix" style (style index 1) on rows 13-15;
# setting .Value above reset each cell back to the default style.
$ws.Range("A16").Copy()
$null = $ws.Range("A13").PasteSpecial(-4122)
$null = $ws.Range("A14").PasteSpecial(-4122)
$null = $ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Selection moves to F7 (from E3).
# ---------------------------------------------------------------------------
$null = $ws.Range("F7").Select()
